$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 0.70407786073899
    3 = 0.0548035619939158
    4 = 2.97911977060175
    5 = 0.270347664275388
    6 = 0.0258834252459316
    7 = 1.81110681072827
    8 = 0.0516840680483948
    9 = 1.43717883140433
    10 = 0.0989006886602981
    11 = 1.97149655343502
    12 = 1.21312402759746
    13 = 1.88810600406152
    14 = 3.3832477881639
    15 = 0.819690893597734
    16 = 2.1980387098307
    17 = 1.09744812064912
    18 = 0.228081808263319
    19 = 13.7720731342522
    20 = 0.160664485015372
    21 = 0.643974643779113
    22 = 0.0653725222608279
    23 = 0.123638652774083
    24 = 10.2082413550787
    25 = 3.9417212396498
    26 = 0.0276872742399872
    27 = 5.0953440634187
    28 = 0.438195356514919
    29 = 1.7113955155636
    30 = 0.650232374672326
    31 = 0.260707163771493
    32 = 0.0747640023716853
    33 = 0.064589655110497
    34 = 3.21093110680525
    35 = 0.774602237130209
    36 = 0.672798825987034
    37 = 0.617223301277396
    38 = 0.343576896536738
    39 = 0.0988419284664624
    40 = 0.134797821354429
    41 = 0.29132237466126
    42 = 9.86440193772741
    43 = 0.138619155658683
    44 = 0.506594327255645
    45 = 0.354025920959256
    46 = 0.0711042623341825
    47 = 0.816026767596802
    48 = 0.039158240401769
    49 = 1.5315850689044
    50 = 0.013650264423317
    51 = 1.10691228621821
    52 = 1.98487486740324
    53 = 0.572731934767874
    54 = 1.02316919575
    55 = 1.644746456904
    56 = 1.8094086365784
    57 = 0.380743171964538
    58 = 0.0593878663199646
    59 = 4.06677526382584
    60 = 1.02012829742575
    61 = 2.36308918374554
    62 = 0.0363768212229457
    63 = 0.260012215827832
    64 = 0.137387917745913
    65 = 0.702072981782492
    66 = 3.55441805002431
    67 = 8.17730278933716
    68 = 0.206108191076833
    69 = 0.521968908394629
    70 = 0.315388615837653
    71 = 1.21958786980526
    72 = 1.03437440382043
    73 = 0.137034341685747
    74 = 6.1097534378829
    75 = 7.29814346652906
    76 = 0.0030762681673735
    77 = 2.19255623489079
    78 = 0.420791042408278
    79 = 0.0222031222618622
    80 = 0.0375574322300205
    81 = 0.949922245841844
    82 = 1.00630953842921
    83 = 0.629510324528561
    84 = 1.05467975514263
    85 = 0.153706983387492
    86 = 1.69167375852436
    87 = 0.00683553092148555
    88 = 0.377780510058461
    89 = 2.484313979195
    90 = 1.08428027488441
    91 = 1.97064061646523
    92 = 1.11351997915374
    93 = 1.07877299890844
    94 = 2.59059663310551
    95 = 1.28648514199607
    96 = 1.53321394040202
    97 = 3.96450783702067
    98 = 5.35385891395565
    99 = 2.57458468836714
    100 = 4.696697865759
    101 = 2.88765603668882
    102 = 0.446642600756921
    103 = 24.775200728038
    104 = 0.773826339611902
    105 = 5.55179788934782
    106 = 0.453718286809018
    107 = 1.42911039196314
    108 = 9.06097578486476
    109 = 0.953618352921304
    110 = 0.665467390942287
    111 = 0.568579162139751
    112 = 1.55793048059032
    113 = 1.0901145368176
    114 = 0.92847757431598
    115 = 1.07580104352365
    116 = 0.052631977764932
    117 = 0.360188852588105
    118 = 0.0501279187512685
    119 = 0.00681329820144433
    120 = 0.37580033453623
    121 = 0.408177737914038
    122 = 1.08623093661785
    123 = 0.0606456211849658
    124 = 1.35662128346383
    125 = 0.525822815244544
    126 = 3.4690781587995
    127 = 0.445700403043463
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row]
}

Write-Host "Done updating $($values.Count) cells"